$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.622051477432251
$ws.Range("B1").Value = 1.92631459236145
$ws.Range("C1").Value = 2.046913146972656
$ws.Range("D1").Value = 2.380285978317261
$ws.Range("E1").Value = 3.188485860824585
